# Update cryptocurrency price (D) and volume change (E) columns
# to reflect the latest scraped values from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.494.93'
$ws.Range("E2").Value = '  +5.99%  '

$ws.Range("D3").Value = '3.559.29'
$ws.Range("E3").Value = '  +2.51%  '

$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").Value = '''417.65'
$ws.Range("E5").Value = '  +0.56%  '

$ws.Range("D6").Value = '''129.99'
$ws.Range("E6").Value = '  +0.36%  '

$ws.Range("D7").Value = '''0.654'
$ws.Range("E7").Value = '  +4.43%  '

$ws.Range("D8").Value = '3.553.46'
$ws.Range("E8").Value = '  +2.63%  '

$ws.Range("D10").Value = '''0.777'
$ws.Range("E10").Value = '  +6.30%  '

$ws.Range("D11").Value = '''0.182'
$ws.Range("E11").Value = '  +29.16%  '

$ws.Range("D12").Value = '''0.0000344'
$ws.Range("E12").Value = '  +56.64%  '

$ws.Range("D13").Value = '''43.05'
$ws.Range("E13").Value = '  +0.53%  '

$ws.Range("D14").Value = '''10.06'
$ws.Range("E14").Value = '  +5.36%  '

$ws.Range("D15").Value = '4.114.25'
$ws.Range("E15").Value = '  +2.51%  '

$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("D17").Value = '''20.37'
$ws.Range("E17").Value = '  -1.06%  '

$ws.Range("D18").Value = '3.617.56'
$ws.Range("E18").Value = '  +5.25%  '

$ws.Range("D19").Value = '''1.12'
$ws.Range("E19").Value = '  +4.39%  '

$ws.Range("D20").Value = '''12.53'
$ws.Range("E20").Value = '  -2.67%  '

$ws.Range("D21").Value = '66.333.79'
$ws.Range("E21").Value = '  +5.72%  '

$ws.Range("D22").Value = '''447.92'
$ws.Range("E22").Value = '  -5.45%  '

$ws.Range("D23").Value = '''89.96'
$ws.Range("E23").Value = '  -1.16%  '

$ws.Range("E24").Value = '  -2.75%  '

$ws.Range("D25").Value = '''13.11'
$ws.Range("E25").Value = '  -2.11%  '

$ws.Range("D26").Value = '''3.37'
$ws.Range("E26").Value = '  +1.76%  '

$ws.Range("D27").Value = '''9.95'
$ws.Range("E27").Value = '  -5.40%  '

$ws.Range("D28").Value = '''34.53'
$ws.Range("E28").Value = '  +3.10%  '

$ws.Range("E29").Value = '  +0.53%  '

$ws.Range("E30").Value = '  +5.16%  '

$ws.Range("D31").Value = '''12.45'
$ws.Range("E31").Value = '  +3.79%  '

$ws.Range("D32").Value = '''0.117'
$ws.Range("E32").Value = '  +4.30%  '

$ws.Range("D33").Value = '''7.26'
$ws.Range("E33").Value = '  -4.88%  '

$ws.Range("D34").Value = '''0.159'
$ws.Range("E34").Value = '  -4.40%  '

$ws.Range("D35").Value = '''0.998'
$ws.Range("E35").Value = '  -0.29%  '

$ws.Range("D36").Value = '''39.00'
$ws.Range("E36").Value = '  -4.45%  '

$ws.Range("D37").Value = '0.0₃0816'
$ws.Range("E37").Value = '  +48.41%  '

$ws.Range("D38").Value = '''56.82'
$ws.Range("E38").Value = '  -2.74%  '

$ws.Range("D39").Value = '''0.0497'
$ws.Range("E39").Value = '  +1.15%  '

$ws.Range("E40").Value = '  +8.95%  '

$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("D42").Value = '''2.78'
$ws.Range("E42").Value = '  +3.18%  '

$ws.Range("D43").Value = '''2.99'
$ws.Range("E43").Value = '  -1.48%  '

$ws.Range("D44").Value = '''148.82'
$ws.Range("E44").Value = '  +2.45%  '

$ws.Range("E45").Value = '  +1.09%  '

$ws.Range("D46").Value = '''3.24'
$ws.Range("E46").Value = '  -3.61%  '

$ws.Range("D47").Value = '''0.308'
$ws.Range("E47").Value = '  -5.22%  '

$ws.Range("D48").Value = '''1.98'
$ws.Range("E48").Value = '  -4.96%  '

$ws.Range("D49").Value = '''2.41'
$ws.Range("E49").Value = '  +0.43%  '

$ws.Range("D50").Value = '''15.55'
$ws.Range("E50").Value = '  -5.37%  '

$ws.Range("E51").Value = '  +10.11%  '
